# Update "想去人数" (interest count) values in F column across sheets,
# matching the generated-output refresh recorded in the commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 107
$ws1.Range("F5").Value  = 287
$ws1.Range("F7").Value  = 1150
$ws1.Range("F17").Value = 495
$ws1.Range("F24").Value = 1937
$ws1.Range("F25").Value = 2482
$ws1.Range("F26").Value = 1255
$ws1.Range("F27").Value = 54
$ws1.Range("F28").Value = 185
$ws1.Range("F29").Value = 355
$ws1.Range("F30").Value = 775
$ws1.Range("F32").Value = 944
$ws1.Range("F33").Value = 115
$ws1.Range("F35").Value = 743
$ws1.Range("F36").Value = 348
$ws1.Range("F38").Value = 711

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 326

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 107
$ws4.Range("F8").Value  = 287
$ws4.Range("F12").Value = 1150
$ws4.Range("F21").Value = 1312
$ws4.Range("F22").Value = 495
$ws4.Range("F27").Value = 2482
$ws4.Range("F29").Value = 1255
$ws4.Range("F30").Value = 54
$ws4.Range("F34").Value = 185
$ws4.Range("F35").Value = 355
$ws4.Range("F36").Value = 775
$ws4.Range("F40").Value = 944
$ws4.Range("F41").Value = 743
$ws4.Range("F42").Value = 348
$ws4.Range("F44").Value = 711

$wb.Save()
